$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9815.949670760207
$ws.Range("D2").Value = 11080.90131642112
$ws.Range("C3").Value = 9826.007148459063
$ws.Range("D3").Value = 11088.27459333561
$ws.Range("C4").Value = 9824.061514558211
$ws.Range("D4").Value = 11086.79065455147
$ws.Range("C5").Value = 9821.019556162446
$ws.Range("D5").Value = 11085.24026376805
$ws.Range("C6").Value = 9819.171638125825
$ws.Range("D6").Value = 11085.39793521137
$ws.Range("C7").Value = 9825.431397126595
$ws.Range("D7").Value = 11091.3507409345
$ws.Range("C8").Value = 9818.655651756611
$ws.Range("D8").Value = 11084.70910010637
$ws.Range("C9").Value = 9821.626720457316
$ws.Range("D9").Value = 11086.23163825032
$ws.Range("C10").Value = 9823.755602750665
$ws.Range("D10").Value = 11088.18033708643
$ws.Range("C11").Value = 9837.442405142605
$ws.Range("D11").Value = 11090.93397972253
$ws.Range("C12").Value = 9819.942370538216
$ws.Range("D12").Value = 11093.69777182456
$ws.Range("C13").Value = 9831.078282645598
$ws.Range("D13").Value = 11090.60443266814
$ws.Range("C14").Value = 9821.558412649962
$ws.Range("D14").Value = 11086.9784376552
$ws.Range("C15").Value = 9825.004388552461
$ws.Range("D15").Value = 11088.92345206815
$ws.Range("C16").Value = 9857.385172610882
$ws.Range("D16").Value = 11100.35375619341
$ws.Range("C17").Value = 9837.962805443494
$ws.Range("D17").Value = 11103.57271551856
$ws.Range("C18").Value = 9826.299665939296
$ws.Range("D18").Value = 11093.07073762144
$ws.Range("C19").Value = 9829.125671183921
$ws.Range("D19").Value = 11090.53314665118
$ws.Range("C20").Value = 9853.849580392331
$ws.Range("D20").Value = 11125.53790092218
$ws.Range("C21").Value = 9844.532197120481
$ws.Range("D21").Value = 11095.17896617037
$ws.Range("C22").Value = 9836.316324579329
$ws.Range("D22").Value = 11101.34295493952
$ws.Range("C23").Value = 9822.107696063355
$ws.Range("D23").Value = 11094.96211378789
$ws.Range("C24").Value = 9862.173181859631
$ws.Range("D24").Value = 11128.39182610468
$ws.Range("C25").Value = 9825.184874868022
$ws.Range("D25").Value = 11087.08251758934
$ws.Range("C26").Value = 9863.893784801983
$ws.Range("D26").Value = 11111.86314032072
$ws.Range("C27").Value = 9898.945484241593
$ws.Range("D27").Value = 11164.21496215561
